# The deck's theme (ppt/theme/theme2.xml, applied via the slide master) is
# switched from the "Integral" / Red Violet colour scheme back to the
# default "Office Theme" / Office colour scheme. (The companion theme part
# used only by the notes master picks up the old "Integral" colours, but
# that part isn't reachable through the PowerPoint object model, so we
# reproduce the visible/reachable half of the swap: the slide master's
# applied theme colours.)

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$scheme = $m.Theme.ThemeColorScheme

# Office Theme colour scheme (target state), in the fixed
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order (indices 1-12).
$officeColors = @(
    (RGBVal 0x00 0x00 0x00),   # 1  dk1
    (RGBVal 0xFF 0xFF 0xFF),   # 2  lt1
    (RGBVal 0x44 0x54 0x6A),   # 3  dk2
    (RGBVal 0xE7 0xE6 0xE6),   # 4  lt2
    (RGBVal 0x5B 0x9B 0xD5),   # 5  accent1
    (RGBVal 0xED 0x7D 0x31),   # 6  accent2
    (RGBVal 0xA5 0xA5 0xA5),   # 7  accent3
    (RGBVal 0xFF 0xC0 0x00),   # 8  accent4
    (RGBVal 0x44 0x72 0xC4),   # 9  accent5
    (RGBVal 0x70 0xAD 0x47),   # 10 accent6
    (RGBVal 0x05 0x63 0xC1),   # 11 hlink
    (RGBVal 0x95 0x4F 0x72)    # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
